$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 12.03168596926108
$ws.Range("D2").Value = 2.904077580570493
$ws.Range("E2").Value = 11.56357016728354
$ws.Range("F2").Value = 66.46222461289308
$ws.Range("G2").Value = 3.859358581380337
$ws.Range("J2").Value = 11.67077120899831
$ws.Range("K2").Value = 30.71353641962546
$ws.Range("L2").Value = 8.270434650084781
$ws.Range("M2").Value = 27.51952848176111
$ws.Range("N2").Value = 24.7244758906123
$ws.Range("C3").Value = 12.0300295481218
$ws.Range("D3").Value = 2.918247984212446
$ws.Range("E3").Value = 11.58634980321889
$ws.Range("F3").Value = 66.44751764270458
$ws.Range("G3").Value = 3.864102976157008
$ws.Range("J3").Value = 11.69770104946853
$ws.Range("K3").Value = 30.5620569267584
$ws.Range("L3").Value = 8.276553785760727
$ws.Range("M3").Value = 27.49003223578289
$ws.Range("N3").Value = 24.79068447089651
$ws.Range("C4").Value = 12.03136910159975
$ws.Range("D4").Value = 2.927352279820081
$ws.Range("E4").Value = 11.60147443058182
$ws.Range("F4").Value = 66.45274908543676
$ws.Range("G4").Value = 3.867164866412808
$ws.Range("J4").Value = 11.71534830569874
$ws.Range("K4").Value = 30.47660953590291
$ws.Range("L4").Value = 8.280539600124758
$ws.Range("M4").Value = 27.47774069103881
$ws.Range("N4").Value = 24.83339750457418
$ws.Range("C5").Value = 12.03250757667621
$ws.Range("D5").Value = 2.931163584074031
$ws.Range("E5").Value = 11.60792448270976
$ws.Range("F5").Value = 66.45845961291427
$ws.Range("G5").Value = 3.868450187912197
$ws.Range("J5").Value = 11.72281989401679
$ws.Range("K5").Value = 30.44371441908522
$ws.Range("L5").Value = 8.282221516609113
$ws.Range("M5").Value = 27.47419754107639
$ws.Range("N5").Value = 24.85132226543892
$ws.Range("C6").Value = 12.03273240588538
$ws.Range("D6").Value = 2.931802555004288
$ws.Range("E6").Value = 11.60901283696563
$ws.Range("F6").Value = 66.45962372786208
$ws.Range("G6").Value = 3.868665888668062
$ws.Range("J6").Value = 11.72407748331712
$ws.Range("K6").Value = 30.43836914269263
$ws.Range("L6").Value = 8.282504285417859
$ws.Range("M6").Value = 27.47369778421371
$ws.Range("N6").Value = 24.85433001585317
$ws.Range("C7").Value = 12.03138205624531
$ws.Range("D7").Value = 2.92740327085269
$ws.Range("E7").Value = 11.60156025692355
$ws.Range("F7").Value = 66.4528116209144
$ws.Range("G7").Value = 3.867182048354507
$ws.Range("J7").Value = 11.71544793498931
$ws.Range("K7").Value = 30.47615807543935
$ws.Range("L7").Value = 8.280562049326164
$ws.Range("M7").Value = 27.47768696953299
$ws.Range("N7").Value = 24.83363714242247
$ws.Range("C8").Value = 12.03062580394368
$ws.Range("D8").Value = 2.908879604870329
$ws.Range("E8").Value = 11.57118874232892
$ws.Range("F8").Value = 66.4541896702564
$ws.Range("G8").Value = 3.8609636568892
$ws.Range("J8").Value = 11.67982608160458
$ws.Range("K8").Value = 30.6597520954478
$ws.Range("L8").Value = 8.272497176750626
$ws.Range("M8").Value = 27.50815213604857
$ws.Range("N8").Value = 24.74687731695053
$ws.Range("C9").Value = 12.04783150892767
$ws.Range("D9").Value = 2.884444526798223
$ws.Range("E9").Value = 11.52063539187909
$ws.Range("F9").Value = 66.57031601690255
$ws.Range("G9").Value = 3.849943044671971
$ws.Range("J9").Value = 11.61877394528605
$ws.Range("K9").Value = 31.07859689671301
$ws.Range("L9").Value = 8.25848843378585
$ws.Range("M9").Value = 27.61392195578066
$ws.Range("N9").Value = 24.59306061955077
$ws.Range("C10").Value = 12.07183847743951
$ws.Range("D10").Value = 2.909764684265369
$ws.Range("E10").Value = 11.48895235577394
$ws.Range("F10").Value = 66.72503830568924
$ws.Range("G10").Value = 3.842551700278513
$ws.Range("J10").Value = 11.57925378376475
$ws.Range("K10").Value = 31.42042312887799
$ws.Range("L10").Value = 8.249286827952748
$ws.Range("M10").Value = 27.71944061532322
$ws.Range("N10").Value = 24.48995451656022
$ws.Range("C11").Value = 12.0852162898052
$ws.Range("D11").Value = 2.921077991876084
$ws.Range("E11").Value = 11.475717617664
$ws.Range("F11").Value = 66.81051042518901
$ws.Range("G11").Value = 3.839340240212544
$ws.Range("J11").Value = 11.56242695157824
$ws.Range("K11").Value = 31.5829092846173
$ws.Range("L11").Value = 8.245335348894606
$ws.Range("M11").Value = 27.77340886371552
$ws.Range("N11").Value = 24.44518993629972
$ws.Range("C12").Value = 12.09063388998308
$ws.Range("D12").Value = 2.925332992417013
$ws.Range("E12").Value = 11.47087485035055
$ws.Range("F12").Value = 66.84504343836748
$ws.Range("G12").Value = 3.838145677435276
$ws.Range("J12").Value = 11.55622011021612
$ws.Range("K12").Value = 31.64540434541426
$ws.Range("L12").Value = 8.243872557610542
$ws.Range("M12").Value = 27.7946953725359
$ws.Range("N12").Value = 24.42854574907623
$ws.Range("C13").Value = 12.08945149878552
$ws.Range("D13").Value = 2.924417890451896
$ws.Range("E13").Value = 11.47191032053408
$ws.Range("F13").Value = 66.83750984834059
$ws.Range("G13").Value = 3.838401991987391
$ws.Range("J13").Value = 11.55754952746652
$ws.Range("K13").Value = 31.63190264578509
$ws.Range("L13").Value = 8.244186106280287
$ws.Range("M13").Value = 27.79007328282939
$ws.Range("N13").Value = 24.43211671606934
$ws.Range("C14").Value = 12.08565496033515
$ws.Range("D14").Value = 2.921428634497587
$ws.Range("E14").Value = 11.4753158170772
$ws.Range("F14").Value = 66.81330808068964
$ws.Range("G14").Value = 3.83924153183922
$ws.Range("J14").Value = 11.56191300417634
$ws.Range("K14").Value = 31.5880316833716
$ws.Range("L14").Value = 8.24521433270465
$ws.Range("M14").Value = 27.77514313457542
$ws.Range("N14").Value = 24.44381445458926
$ws.Range("C15").Value = 12.08337522209914
$ws.Range("D15").Value = 2.919593849714559
$ws.Range("E15").Value = 11.47742376954143
$ws.Range("F15").Value = 66.79876582401052
$ws.Range("G15").Value = 3.83975857585463
$ws.Range("J15").Value = 11.56460725044868
$ws.Range("K15").Value = 31.56128392321362
$ws.Range("L15").Value = 8.245848515452368
$ws.Range("M15").Value = 27.76610841089102
$ws.Range("N15").Value = 24.45101965125989
$ws.Range("C16").Value = 12.07101351640017
$ws.Range("D16").Value = 2.909021240321445
$ws.Range("E16").Value = 11.48984094278237
$ws.Range("F16").Value = 66.71975589852761
$ws.Range("G16").Value = 3.842764600257607
$ws.Range("J16").Value = 11.58037658437868
$ws.Range("K16").Value = 31.40994153219332
$ws.Range("L16").Value = 8.249549769380424
$ws.Range("M16").Value = 27.71603306299876
$ws.Range("N16").Value = 24.49292299379589
$ws.Range("C17").Value = 12.06405825041397
$ws.Range("D17").Value = 2.902483055383287
$ws.Range("E17").Value = 11.49775987157276
$ws.Range("F17").Value = 66.67514878319311
$ws.Range("G17").Value = 3.844647239831963
$ws.Range("J17").Value = 11.59034508642516
$ws.Range("K17").Value = 31.31885991939509
$ws.Range("L17").Value = 8.251880287610456
$ws.Range("M17").Value = 27.68683631013688
$ws.Range("N17").Value = 24.51917687811694
$ws.Range("C18").Value = 12.06028913659435
$ws.Range("D18").Value = 2.898703150665121
$ws.Range("E18").Value = 11.50242553939201
$ws.Range("F18").Value = 66.6509124393876
$ws.Range("G18").Value = 3.845744296762486
$ws.Range("J18").Value = 11.59618707495457
$ws.Range("K18").Value = 31.26713240468282
$ws.Range("L18").Value = 8.253242809963346
$ws.Range("M18").Value = 27.67060533133375
$ws.Range("N18").Value = 24.53447876270511
$ws.Range("C19").Value = 12.05905276119995
$ws.Range("D19").Value = 2.897420011466875
$ws.Range("E19").Value = 11.50402431717317
$ws.Range("F19").Value = 66.64295045952575
$ws.Range("G19").Value = 3.846118187014262
$ws.Range("J19").Value = 11.59818369713645
$ws.Range("K19").Value = 31.24973293027815
$ws.Range("L19").Value = 8.25370793190015
$ws.Range("M19").Value = 27.66520659612378
$ws.Range("N19").Value = 24.53969431121183
$ws.Range("C20").Value = 12.06477471488588
$ws.Range("D20").Value = 2.903181052539793
$ws.Range("E20").Value = 11.49690541271728
$ws.Range("F20").Value = 66.67975028726806
$ws.Range("G20").Value = 3.844445359801526
$ws.Range("J20").Value = 11.58927270950428
$ws.Range("K20").Value = 31.32848766428603
$ws.Range("L20").Value = 8.25162991694431
$ws.Range("M20").Value = 27.6898862227677
$ws.Range("N20").Value = 24.51636127179289
$ws.Range("C21").Value = 12.08676056445448
$ws.Range("D21").Value = 2.922307438359527
$ws.Range("E21").Value = 11.47431095883254
$ws.Range("F21").Value = 66.82035796119787
$ws.Range("G21").Value = 3.838994355049981
$ws.Range("J21").Value = 11.56062686735073
$ws.Range("K21").Value = 31.60089178224993
$ws.Range("L21").Value = 8.24491140842245
$ws.Range("M21").Value = 27.77950548487689
$ws.Range("N21").Value = 24.44037021113067
$ws.Range("C22").Value = 12.10317869491049
$ws.Range("D22").Value = 2.9346381700051
$ws.Range("E22").Value = 11.46052868218637
$ws.Range("F22").Value = 66.9248787620528
$ws.Range("G22").Value = 3.83555733311013
$ws.Range("J22").Value = 11.54286738484905
$ws.Range("K22").Value = 31.78452917379271
$ws.Range("L22").Value = 8.240715947735353
$ws.Range("M22").Value = 27.84302627349031
$ws.Range("N22").Value = 24.39249603069897
$ws.Range("C23").Value = 12.09422913768536
$ws.Range("D23").Value = 2.928072411085709
$ws.Range("E23").Value = 11.46779461125877
$ws.Range("F23").Value = 66.86794036392681
$ws.Range("G23").Value = 3.837380300950661
$ws.Range("J23").Value = 11.55225804387596
$ws.Range("K23").Value = 31.68601912263065
$ws.Range("L23").Value = 8.242937309106177
$ws.Range("M23").Value = 27.80867411295315
$ws.Range("N23").Value = 24.41788368103227
$ws.Range("C24").Value = 12.06445008614915
$ws.Range("D24").Value = 2.90286555351098
$ws.Range("E24").Value = 11.49729136165432
$ws.Range("F24").Value = 66.6776655587169
$ws.Range("G24").Value = 3.844536583957147
$ws.Range("J24").Value = 11.58975718545125
$ws.Range("K24").Value = 31.32413297521868
$ws.Range("L24").Value = 8.251743038878931
$ws.Range("M24").Value = 27.68850562913061
$ws.Range("N24").Value = 24.51763355895863
$ws.Range("C25").Value = 12.04117781078551
$ws.Range("D25").Value = 2.884380683298758
$ws.Range("E25").Value = 11.53335059471315
$ws.Range("F25").Value = 66.52672330440045
$ws.Range("G25").Value = 3.852799799335243
$ws.Range("J25").Value = 11.63435108761463
$ws.Range("K25").Value = 30.95916278267711
$ws.Range("L25").Value = 8.262085883133265
$ws.Range("M25").Value = 27.58040268922352
$ws.Range("N25").Value = 24.63292985880293
